$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("43267", "Le château des princesses et les animaux royaux", "787", "Disney™", "https://www.lego.com/cdn/cs/set/assets/bltc5e58398c856ec9b/43267_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1", "https://amzn.eu/d/0IEMO6h", "https://www.lego.com/fr-fr/product/43267", "https://www.auchan.fr/lego-disney-43267-la-chateau-des-princesses-et-les-animaux-royaux/pr-C1839464", "https://www.e.leclerc/fp/lego-disney-princess-le-chateau-des-princesses-et-les-animaux-royaux-43267-5702017814568", "https://www.carrefour.fr/p/lego-le-chateau-des-princesses-et-les-animaux-royaux-43267-lego-5702017814568", ""),
    @("10342", "Le joli bouquet de fleurs roses", "749", "The Botanical Collection", "https://www.lego.com/cdn/cs/set/assets/blt1bcb187351a12753/10342_WEB_PRI_NOBG.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1", "https://amzn.eu/d/dJwTE0D", "https://www.lego.com/fr-fr/product/10342", "https://www.auchan.fr/lego-botanicals-10342-le-joli-bouquet-de-fleurs-roses/pr-C1821871", "https://www.e.leclerc/fp/lego-botanique-10342-le-joli-bouquet-de-fleurs-roses-set-pour-adultes-5702017812540", "https://www.carrefour.fr/p/lego-le-joli-bouquet-de-fleurs-roses-10342-lego-5702017812540", ""),
    @("31167", "Le manoir hanté", "736", "Creator 3-en-1", "https://www.lego.com/cdn/cs/set/assets/blt61735fb0922b35cc/31167_Prod.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1", "https://amzn.eu/d/9geJ4zv", "https://www.lego.com/fr-fr/product/31167", "https://www.auchan.fr/lego-creator-31167-le-manoir-hante/pr-C1844715", "https://www.e.leclerc/fp/lego-creator-le-manoir-hante-jouet-de-construction-31167-5702017822426", "", ""),
    @("60367", "L’avion de ligne", "913", "City", "https://www.lego.com/cdn/cs/set/assets/blta636ffccc87fdf4b/60367.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1", "https://amzn.eu/d/1x42GYW", "https://www.lego.com/fr-fr/product/60367", "https://www.auchan.fr/lego-lego-city-60367-l-avion-de-ligne-jouet-avec-bus-de-l-aeroport-remorqueur-camion-a-bagages-et-9-minifigurines/pr-C1718310", "https://www.e.leclerc/fp/l-avion-de-ligne-5702017416274?srsltid=AfmBOorK1WUD-ZnAbLCczVAVBtBljFx-AHHh7iXqiXiScD6CX2JFzn2B", "https://www.carrefour.fr/p/lego-city-l-avion-de-ligne-60367-lego-5702017416274", "")
)

$startRow = 18
$endRow = $startRow + $rows.Length - 1

# Force the new block to be treated as text (so purely-numeric strings like
# "43267" stay strings, matching the source data) without leaving a lasting
# number-format override on the cells.
$fullRange = $ws.Range("A" + $startRow + ":K" + $endRow)
$fullRange.NumberFormat = "@"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le 11; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}

# Revert the block back to the default "Normal" style so the cells don't
# carry an explicit style index, matching the rest of the sheet.
$fullRange.Style = "Normal"
